$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.087.29'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '2.314.02'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '302.11'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = '99.14'
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '0.116'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = '17.97'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").Value = '6.94'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '2.674.82'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '2.313.87'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '43.018.22'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '13.48'
$ws.Range("E19").Value = '  +6.79%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '68.06'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '240.67'
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").Value = '24.97'
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").Value = '168.31'
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").Value = "'9.20"
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '33.53'
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").Value = '5.24'
$ws.Range("E32").Value = '  +4.60%  '
$ws.Range("D33").Value = '4.93'
$ws.Range("E33").Value = '  +7.08%  '
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +8.21%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = '0.0694'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '1.998.29'
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("E44").Value = '  -6.17%  '
$ws.Range("E45").Value = '  -1.18%  '
$ws.Range("D46").Value = '17.45'
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '54.95'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").Value = '74.92'
$ws.Range("E49").Value = '  +6.75%  '
$ws.Range("D50").Value = '2.541.21'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("E51").Value = '  +1.14%  '
